$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1242
$ws.Range("I100").Value = 1255.9
$ws.Range("J100").Value = 1195.6666
$ws.Range("K100").Value = 1255.9
$ws.Range("L100").Value = 1195.6666
$ws.Range("M100").Value = -714.9000000000001
$ws.Range("N100").Value = -2277.6666
$ws.Range("H121").Value = 513.6667
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 513.6667
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 1541.0001
$ws.Range("M121").Value = $null
$ws.Range("N121").Value = -5035.0001
$ws.Range("H129").Value = 6732.8823
$ws.Range("J129").Value = 1090.7
$ws.Range("L129").Value = 3272.1
$ws.Range("N129").Value = -13272.1
$ws.Range("H136").Value = 68000
$ws.Range("J136").Value = 68000
$ws.Range("L136").Value = 68000
$ws.Range("N136").Value = -78200
$ws.Range("H141").Value = 4900
$ws.Range("I141").Value = 4950
$ws.Range("J141").Value = 4600
$ws.Range("K141").Value = 14850
$ws.Range("L141").Value = 13800
$ws.Range("M141").Value = -9670
$ws.Range("N141").Value = -24160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = $null
$ws.Range("H32").Value = 14258.69
$ws.Range("I32").Value = 12339.147
$ws.Range("J32").Value = 28335.334
$ws.Range("K32").Value = 12339.147
$ws.Range("L32").Value = 28335.334
$ws.Range("M32").Value = -12052.147
$ws.Range("N32").Value = -28909.334
$ws.Range("H45").Value = 78319.54
$ws.Range("I45").Value = 143758.58
$ws.Range("K45").Value = 143758.58
$ws.Range("M45").Value = -143381.58
$ws.Range("H74").Value = 1488.0834
$ws.Range("I74").Value = 1436
$ws.Range("K74").Value = 1436
$ws.Range("M74").Value = -562
$ws.Range("H77").Value = 1488.0834
$ws.Range("I77").Value = 1436
$ws.Range("K77").Value = 7180
$ws.Range("M77").Value = -2812
$ws.Range("H97").Value = 44583.78
$ws.Range("I97").Value = 50915.25
$ws.Range("J97").Value = 2374
$ws.Range("K97").Value = 50915.25
$ws.Range("L97").Value = 2374
$ws.Range("M97").Value = -50419.25
$ws.Range("N97").Value = -3366
$ws.Range("H122").Value = 2041.75
$ws.Range("I122").Value = 1750.2941
$ws.Range("J122").Value = 3693.3333
$ws.Range("K122").Value = 5250.8823
$ws.Range("L122").Value = 11079.9999
$ws.Range("M122").Value = -2800.8823
$ws.Range("N122").Value = -15979.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1826.0555
$ws.Range("I99").Value = 1559.9166
$ws.Range("K99").Value = 1559.9166
$ws.Range("M99").Value = -61.91660000000002
$ws.Range("H105").Value = 107078.31
$ws.Range("I105").Value = 85274.086
$ws.Range("J105").Value = 144457
$ws.Range("K105").Value = 85274.086
$ws.Range("L105").Value = 144457
$ws.Range("M105").Value = -83527.086
$ws.Range("N105").Value = -147951
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H134").Value = 2873.0376
$ws.Range("I134").Value = 2782.1304
$ws.Range("K134").Value = 8346.3912
$ws.Range("M134").Value = -5811.3912
$ws.Range("H139").Value = 62000
$ws.Range("J139").Value = 62000
$ws.Range("L139").Value = 62000
$ws.Range("N139").Value = -72280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33856.176
$ws.Range("I31").Value = 1415.3334
$ws.Range("J31").Value = 94682.75
$ws.Range("K31").Value = 1415.3334
$ws.Range("L31").Value = 94682.75
$ws.Range("M31").Value = -1120.3334
$ws.Range("N31").Value = -95272.75
$ws.Range("H34").Value = 33856.176
$ws.Range("I34").Value = 1415.3334
$ws.Range("J34").Value = 94682.75
$ws.Range("K34").Value = 1415.3334
$ws.Range("L34").Value = 94682.75
$ws.Range("M34").Value = -1213.3334
$ws.Range("N34").Value = -95086.75
$ws.Range("H105").Value = 1552.9
$ws.Range("I105").Value = 1593.625
$ws.Range("J105").Value = 1390
$ws.Range("K105").Value = 1593.625
$ws.Range("L105").Value = 1390
$ws.Range("M105").Value = 153.375
$ws.Range("N105").Value = -4884
$ws.Range("H122").Value = 1043.1111
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 1117.6
$ws.Range("K122").Value = 2850
$ws.Range("L122").Value = 3352.8
$ws.Range("M122").Value = -400
$ws.Range("N122").Value = -8252.8
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -60060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 7466.579
$ws.Range("I140").Value = 13320.556
$ws.Range("J140").Value = 2198
$ws.Range("K140").Value = 39961.66800000001
$ws.Range("L140").Value = 6594
$ws.Range("M140").Value = -34781.66800000001
$ws.Range("N140").Value = -16954

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2008.0294
$ws.Range("J43").Value = 5258.125
$ws.Range("L43").Value = 5258.125
$ws.Range("N43").Value = -5560.125
$ws.Range("H46").Value = 12485.714
$ws.Range("J46").Value = 12485.714
$ws.Range("L46").Value = 12485.714
$ws.Range("N46").Value = -12797.714
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = $null
$ws.Range("H80").Value = 77095230
$ws.Range("J80").Value = 3672.2222
$ws.Range("L80").Value = 3672.2222
$ws.Range("N80").Value = -5668.2222
$ws.Range("H83").Value = 77095230
$ws.Range("J83").Value = 3672.2222
$ws.Range("L83").Value = 18361.111
$ws.Range("N83").Value = -28345.111
$ws.Range("H122").Value = 2131.1428
$ws.Range("I122").Value = 1812.9333
$ws.Range("K122").Value = 5438.7999
$ws.Range("M122").Value = -2988.7999
$ws.Range("H126").Value = 3714.5833
$ws.Range("I126").Value = 3846.25
$ws.Range("J126").Value = 3451.25
$ws.Range("K126").Value = 11538.75
$ws.Range("L126").Value = 10353.75
$ws.Range("M126").Value = -9068.75
$ws.Range("N126").Value = -15293.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 335793.34
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 3690
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 3690
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -3962
$ws.Range("H122").Value = 6156.5
$ws.Range("I122").Value = 4502
$ws.Range("K122").Value = 13506
$ws.Range("M122").Value = -11056

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 100001496
$ws.Range("J96").Value = 952
$ws.Range("L96").Value = 952
$ws.Range("N96").Value = -3698
$ws.Range("H119").Value = 29599.666
$ws.Range("J119").Value = 29599.666
$ws.Range("L119").Value = 29599.666
$ws.Range("N119").Value = -39275.666
